# Update Format file Excel
# Insert a new "SO LUONG" (quantity) column before the "GIA NHAP" column
# (current column E), shifting GIA NHAP / GIA BAN / MO TA one column to
# the right, and fill in a quantity value of 20 for each product row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column at E (existing E:G shift to F:H)
$ws.Columns("E:E").Insert()

# Match the new column's width to its neighbours (B:D)
$ws.Columns("E:E").ColumnWidth = $ws.Columns("D:D").ColumnWidth

# Header for the new column
$ws.Range("E1").Value = "SỐ LƯỢNG"

# Quantity values for the existing product rows
$ws.Range("E2").Value = 20
$ws.Range("E3").Value = 20
$ws.Range("E4").Value = 20
$ws.Range("E5").Value = 20

# Restore the (now shifted) current selection
[void]$ws.Range("E14").Select()

Write-Host "Inserted 'SO LUONG' column and updated Sheet1 formatting."
